$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2082111436950147
$ws.Range("C2").Value = 0.5219941348973607
$ws.Range("J2").Value = 0.02932551319648094
$ws.Range("P2").Value = 0.1319648093841642
$ws.Range("S2").Value = 0.1085043988269795
$ws.Range("B3").Value = 0.01092896174863388
$ws.Range("C3").Value = 0.0273224043715847
$ws.Range("J3").Value = 0.0546448087431694
$ws.Range("P3").Value = 0.6775956284153005
$ws.Range("S3").Value = 0.2295081967213115
$ws.Range("O4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.6060606060606061
$ws.Range("S4").Value = 0.3636363636363636
$ws.Range("B6").Value = 0.06161137440758294
$ws.Range("D6").Value = 0.004739336492890996
$ws.Range("E6").Value = 0.004739336492890996
$ws.Range("F6").Value = 0.05687203791469194
$ws.Range("J6").Value = 0.2369668246445498
$ws.Range("O6").Value = 0.04739336492890995
$ws.Range("Q6").Value = 0.1848341232227488
$ws.Range("R6").Value = 0.07582938388625593
$ws.Range("S6").Value = 0.3270142180094787
$ws.Range("B7").Value = 0.1416309012875537
$ws.Range("D7").Value = 0.0128755364806867
$ws.Range("F7").Value = 0.04721030042918455
$ws.Range("J7").Value = 0.1373390557939914
$ws.Range("O7").Value = 0.008583690987124463
$ws.Range("Q7").Value = 0.1630901287553648
$ws.Range("R7").Value = 0.07725321888412018
$ws.Range("S7").Value = 0.4120171673819742
$ws.Range("B8").Value = 0.09808102345415778
$ws.Range("D8").Value = 0.01492537313432836
$ws.Range("E8").Value = 0.002132196162046908
$ws.Range("F8").Value = 0.05543710021321962
$ws.Range("J8").Value = 0.138592750533049
$ws.Range("O8").Value = 0.01705756929637527
$ws.Range("Q8").Value = 0.1428571428571428
$ws.Range("R8").Value = 0.1087420042643923
$ws.Range("S8").Value = 0.4221748400852878
$ws.Range("B9").Value = 0.0880503144654088
$ws.Range("D9").Value = 0.006289308176100629
$ws.Range("F9").Value = 0.06918238993710692
$ws.Range("J9").Value = 0.1069182389937107
$ws.Range("O9").Value = 0.02515723270440252
$ws.Range("Q9").Value = 0.1949685534591195
$ws.Range("R9").Value = 0.119496855345912
$ws.Range("S9").Value = 0.389937106918239
$ws.Range("B10").Value = 0.1196063588190765
$ws.Range("D10").Value = 0.0174110522331567
$ws.Range("E10").Value = 0.000757002271006813
$ws.Range("F10").Value = 0.06888720666161999
$ws.Range("J10").Value = 0.1362604087812263
$ws.Range("O10").Value = 0.01514004542013626
$ws.Range("Q10").Value = 0.1786525359576079
$ws.Range("R10").Value = 0.08629825889477669
$ws.Range("S10").Value = 0.3769871309613929
$ws.Range("G11").Value = 0.145945945945946
$ws.Range("J11").Value = 0.0972972972972973
$ws.Range("K11").Value = 0.2081081081081081
$ws.Range("L11").Value = 0.5405405405405406
$ws.Range("S11").Value = 0.008108108108108109
$ws.Range("G12").Value = 0.7255813953488373
$ws.Range("J12").Value = 0.1906976744186047
$ws.Range("L12").Value = 0.04186046511627907
$ws.Range("S12").Value = 0.04186046511627907
$ws.Range("G13").Value = 0.8
$ws.Range("J13").Value = 0.1714285714285714
$ws.Range("S13").Value = 0.02857142857142857
$ws.Range("F15").Value = 0.0326530612244898
$ws.Range("H15").Value = 0.1836734693877551
$ws.Range("I15").Value = 0.05714285714285714
$ws.Range("J15").Value = 0.3714285714285714
$ws.Range("K15").Value = 0.0653061224489796
$ws.Range("M15").Value = 0.004081632653061225
$ws.Range("N15").Value = 0.004081632653061225
$ws.Range("O15").Value = 0.05714285714285714
$ws.Range("S15").Value = 0.2244897959183673
$ws.Range("F16").Value = 0.005405405405405406
$ws.Range("H16").Value = 0.1513513513513514
$ws.Range("I16").Value = 0.06486486486486487
$ws.Range("J16").Value = 0.4540540540540541
$ws.Range("K16").Value = 0.0918918918918919
$ws.Range("M16").Value = 0.01621621621621622
$ws.Range("N16").Value = 0.005405405405405406
$ws.Range("O16").Value = 0.08648648648648649
$ws.Range("S16").Value = 0.1243243243243243
$ws.Range("F17").Value = 0.004938271604938272
$ws.Range("H17").Value = 0.1901234567901235
$ws.Range("I17").Value = 0.05432098765432099
$ws.Range("J17").Value = 0.4444444444444444
$ws.Range("K17").Value = 0.1209876543209877
$ws.Range("M17").Value = 0.01975308641975309
$ws.Range("O17").Value = 0.07901234567901234
$ws.Range("S17").Value = 0.08641975308641975
$ws.Range("F18").Value = 0.01376146788990826
$ws.Range("H18").Value = 0.1880733944954129
$ws.Range("I18").Value = 0.09174311926605505
$ws.Range("J18").Value = 0.4220183486238532
$ws.Range("K18").Value = 0.1009174311926606
$ws.Range("M18").Value = 0.02293577981651376
$ws.Range("O18").Value = 0.05504587155963303
$ws.Range("S18").Value = 0.1055045871559633
$ws.Range("F19").Value = 0.01640625
$ws.Range("H19").Value = 0.21953125
$ws.Range("I19").Value = 0.07109375
$ws.Range("J19").Value = 0.346875
$ws.Range("K19").Value = 0.14296875
$ws.Range("M19").Value = 0.01640625
$ws.Range("N19").Value = 0.0015625
$ws.Range("O19").Value = 0.075
$ws.Range("S19").Value = 0.11015625
